$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: paragraph "*Poder añadir un usuario ... -Numero de likes. "
# Change "...-Numero de likes. " into
# "...-Numero de likes" / " y una lista con sus tweets" / bookmark _GoBack /
# ". "  (four separate pieces, matching the target canonical OOXML run split)
# ---------------------------------------------------------------------------

$boundaryFind = $d.Content
$boundaryFind.Find.Execute("-Numero de tweets")
$boundary = $boundaryFind.End

$endFind = $d.Content
$endFind.Find.Execute(" -Numero de seguidores -Numero de seguidos -Numero de likes. ")
$endPos = $endFind.End

$target1 = $d.Range($boundary, $endPos)

$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve"> -Numero de seguidores -Numero de seguidos -Numero de likes</w:t></w:r><w:r><w:t xml:space="preserve"> y una lista con sus tweets</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Edit 2: paragraph "*Conectar usuarios ..." ending
# "...de los 2 usuario" + bookmark _GoBack + "s" + "."  ->
# "...de los 2 usuarios." as a single run, bookmark removed from here
# (it now lives in paragraph edited above).
# ---------------------------------------------------------------------------

$r2 = $d.Content
$r2.Find.Execute("la afinidad se determina como el valor absoluto de la resta entre las puntuaciones de los 2 usuario")
$target2 = $d.Range($r2.Start, $r2.End + 2)

$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">la afinidad se determina como el valor absoluto de la resta entre las puntuaciones de los 2 usuarios.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target2.InsertXML($xml2)
